$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Row 2 values: Day, then 24 hourly prices, Price_Daily_Avg, Slot_4h_max,
# Slot_4h_price, Slot_2h_frist, Slot_2h_frist_price, Slot_2h_second,
# Slot_2h_second_price, Slot_min_price
$ws.Range("A2").Value = 45888
$ws.Range("B2").Value = 103.06
$ws.Range("C2").Value = 101.12
$ws.Range("D2").Value = 98.09999999999999
$ws.Range("E2").Value = 92.41
$ws.Range("F2").Value = 90.37
$ws.Range("G2").Value = 91.40000000000001
$ws.Range("H2").Value = 99
$ws.Range("I2").Value = 102.82
$ws.Range("J2").Value = 98.09999999999999
$ws.Range("K2").Value = 88.56999999999999
$ws.Range("L2").Value = 69.64
$ws.Range("M2").Value = 43.23
$ws.Range("N2").Value = 30
$ws.Range("O2").Value = 29.14
$ws.Range("P2").Value = 30.91
$ws.Range("Q2").Value = 40.54
$ws.Range("R2").Value = 46.8
$ws.Range("S2").Value = 48.7
$ws.Range("T2").Value = 62.2
$ws.Range("U2").Value = 88.97
$ws.Range("V2").Value = 97.73999999999999
$ws.Range("W2").Value = 106.85
$ws.Range("X2").Value = 106.54
$ws.Range("Y2").Value = 99
$ws.Range("Z2").Value = 77.72
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 102.53
$ws.Range("AC2").Value = "22h-24h"
$ws.Range("AD2").Value = 102.77
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 102.29
$ws.Range("AG2").Value = "10h-18h"
